$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F3: "deltaL" -> "deltaL[mm]"
$ws.Range("F3").Value = "deltaL[mm]"

# Row 4: swap the "T_0:" label and the "1" index between A4 and B4
$ws.Range("A4").Value = "T_0:"
$ws.Range("B4").Value = 1

# Rows 5-30: move the running index from column A to column B
$indexRows = @(
  @{r=5;  v=2},  @{r=6;  v=3},  @{r=7;  v=4},  @{r=8;  v=5},
  @{r=9;  v=6},  @{r=10; v=7},  @{r=11; v=8},  @{r=12; v=9},
  @{r=13; v=10}, @{r=14; v=11}, @{r=15; v=12}, @{r=16; v=13},
  @{r=17; v=14}, @{r=18; v=15}, @{r=19; v=16}, @{r=20; v=17},
  @{r=21; v=18}, @{r=22; v=19}, @{r=23; v=20}, @{r=24; v=21},
  @{r=25; v=22}, @{r=26; v=23}, @{r=27; v=24}, @{r=28; v=25},
  @{r=29; v=26}, @{r=30; v=27}
)
foreach ($row in $indexRows) {
  $ws.Cells.Item($row.r, 1).ClearContents()
  $ws.Cells.Item($row.r, 2).Value = $row.v
}

# New measurement data in C5:F8
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = 0.52
$ws.Range("E5").Value = 1.9
$ws.Range("F5").Value = 15

$ws.Range("C6").Value = 38.5
$ws.Range("F6").Value = 24

$ws.Range("C7").Value = 45.3
$ws.Range("F7").Value = 35

$ws.Range("C8").Value = 51.5
$ws.Range("F8").Value = 44

# Update the selected cell
[void]$ws.Range("C9").Select()
